# Create new keyword tag rows.
# 1. Split the old "tag/question" pair that lived in row 258 into:
#      - row 258 keeps the tag "ตกลง" but gets a brand-new question "โอเคร"
#      - row 259 restores the original tag "ถูกหนึ่ง" with a corrected question
#      - row 260 is a fully new tag/question pair
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A258").Value = "ตกลง"
$ws.Range("B258").Value = "โอเคร"

$ws.Range("A259").Value = "ถูกหนึ่ง"
$ws.Range("B259").Value = "ประเทศไทย เพราะประเทศไทยมีตรัง (ตัง)"

$ws.Range("A260").Value = "ถูกสอง"
$ws.Range("B260").Value = "ไปฉันเพล"

$ws.Range("B260").Select()
